# Apply updated cryptocurrency data (row-for-row) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '66.847.29'
$ws.Range('E2').Value = '  +1.42%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.356.92'
$ws.Range('E3').Value = '  +1.49%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '190.82'
$ws.Range('E5').Value = '  +5.55%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '562.40'
$ws.Range('E6').Value = '  +0.63%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.348.69'
$ws.Range('E8').Value = '  +1.43%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.585'
$ws.Range('E9').Value = '  -1.02%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.184'
$ws.Range('E10').Value = '  -2.76%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.588'
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '47.23'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000272'
$ws.Range('E13').Value = '  +2.35%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '8.74'
$ws.Range('E14').Value = '  +1.88%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.882.89'
$ws.Range('E15').Value = '  +1.45%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '606.75'
$ws.Range('E16').Value = '  -4.73%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '66.665.35'
$ws.Range('E17').Value = '  +1.20%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '18.09'
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.118'
$ws.Range('E19').Value = '  +1.28%  '
$ws.Range('B20').Value = 'WrappedEther'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.346.77'
$ws.Range('E20').Value = '  +1.52%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.14'
$ws.Range('E21').Value = '  -2.80%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.909'
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '18.55'
$ws.Range('E23').Value = '  +4.89%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.07'
$ws.Range('E24').Value = '  +0.45%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '100.90'
$ws.Range('E25').Value = '  -6.40%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '4.03'
$ws.Range('E26').Value = '  +0.29%  '
$ws.Range('E27').Value = '  +1.34%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.77'
$ws.Range('E28').Value = '  +3.05%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.72'
$ws.Range('E29').Value = '  +1.80%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.78'
$ws.Range('E30').Value = '  -0.27%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '30.81'
$ws.Range('E31').Value = '  +0.27%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.82'
$ws.Range('E32').Value = '  +7.65%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.98'
$ws.Range('E33').Value = '  -0.18%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '588.18'
$ws.Range('E34').Value = '  +6.67%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('E36').Value = '  -0.31%  '
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.717.02'
$ws.Range('E37').Value = '  -0.61%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '57.24'
$ws.Range('E38').Value = '  +0.26%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  -0.06%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.63'
$ws.Range('E40').Value = '  +7.28%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '34.02'
$ws.Range('E41').Value = '  +5.86%  '
$ws.Range('B42').Value = 'PEPE'
$ws.Range('C42').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0₃0716'
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.130'
$ws.Range('E43').Value = '  +2.14%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.28'
$ws.Range('E44').Value = '  -5.82%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.72'
$ws.Range('E45').Value = '  -0.89%  '
$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.344'
$ws.Range('E46').Value = '  +0.31%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.40'
$ws.Range('E47').Value = '  +5.17%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0425'
$ws.Range('E48').Value = '  +2.58%  '
$ws.Range('E49').Value = '  +0.34%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.61'
$ws.Range('E50').Value = '  -0.37%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.998'
$ws.Range('E51').Value = '  +0.01%  '
